$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("O4").Value = 1.67
$ws.Range("P4").Value = 2.1
$ws.Range("I7").Value = 3.9
$ws.Range("Q7").Value = 3.4
$ws.Range("R7").Value = 1.33
$ws.Range("W7").Value = 5
$ws.Range("AE7").Value = 23
$ws.Range("AH7").Value = 17
$ws.Range("AI7").Value = 15
$ws.Range("K12").Value = 2.18
$ws.Range("AH12").Value = 24
$ws.Range("AM12").Value = 400
$ws.Range("AR12").Value = 55
$ws.Range("AT12").Value = 2.87
$ws.Range("AU12").Value = 7
$ws.Range("AW12").Value = 5.9
$ws.Range("BB12").Value = 350
$ws.Range("M15").Value = 1.05
$ws.Range("N15").Value = 11
$ws.Range("Q15").Value = 2
$ws.Range("R15").Value = 1.8
$ws.Range("G17").Value = 2.6
$ws.Range("I17").Value = 2.55
$ws.Range("M17").Value = 1.07
$ws.Range("N17").Value = 9
$ws.Range("W17").Value = 8
$ws.Range("Y17").Value = 10
$ws.Range("Z17").Value = 26
$ws.Range("AD17").Value = 6
$ws.Range("AM17").Value = 251
$ws.Range("AR17").Value = 67
$ws.Range("G19").Value = 2.5
$ws.Range("H19").Value = 2.9
$ws.Range("Z19").Value = 23
$ws.Range("AA19").Value = 23
$ws.Range("AC19").Value = 6.5
$ws.Range("M20").Value = 1.07
$ws.Range("N20").Value = 9
$ws.Range("Q20").Value = 2.25
$ws.Range("R20").Value = 1.62
$ws.Range("M22").Value = 1.05
$ws.Range("O22").Value = 1.33
$ws.Range("I23").Value = 2.88
$ws.Range("J23").Value = 3.2
$ws.Range("M23").Value = 1.05
$ws.Range("O23").Value = 1.33
$ws.Range("AK23").Value = 26
$ws.Range("AR23").Value = 67
$ws.Range("M24").Value = 1.04
$ws.Range("N24").Value = 10
$ws.Range("O24").Value = 1.3
$ws.Range("G25").Value = 2.25
$ws.Range("I25").Value = 3.5
$ws.Range("M25").Value = 1.05
$ws.Range("O25").Value = 1.33
$ws.Range("AA25").Value = 19
$ws.Range("AG25").Value = 9.5
$ws.Range("AH25").Value = 17
$ws.Range("AI25").Value = 13
$ws.Range("AP25").Value = 23
$ws.Range("BA25").Value = 101
$ws.Range("J27").Value = 2.8
$ws.Range("L27").Value = 3.65
$ws.Range("N27").Value = 6.9
$ws.Range("O27").Value = 1.33
$ws.Range("P27").Value = 3.05
$ws.Range("Q27").Value = 2
$ws.Range("R27").Value = 1.75
$ws.Range("U27").Value = 1.78
$ws.Range("V27").Value = 1.93
$ws.Range("Y27").Value = 9
$ws.Range("AB27").Value = 29
$ws.Range("AC27").Value = 6.9
$ws.Range("AE27").Value = 14
$ws.Range("AG27").Value = 9
$ws.Range("AH27").Value = 15.5
$ws.Range("AI27").Value = 11
$ws.Range("AK27").Value = 28
$ws.Range("AL27").Value = 37
$ws.Range("AO27").Value = 11.5
$ws.Range("AR27").Value = 75
$ws.Range("AU27").Value = 7
$ws.Range("BA27").Value = 120
$ws.Range("G37").Value = 2.25
$ws.Range("H37").Value = 3.15
$ws.Range("J37").Value = 2.87
$ws.Range("K37").Value = 2.02
$ws.Range("L37").Value = 3.65
$ws.Range("S37").Value = 1.44
$ws.Range("T37").Value = 2.42
$ws.Range("U37").Value = 1.83
$ws.Range("V37").Value = 1.78
$ws.Range("X37").Value = 10.25
$ws.Range("AA37").Value = 20
$ws.Range("AC37").Value = 8
$ws.Range("AD37").Value = 6.1
$ws.Range("AE37").Value = 15.5
$ws.Range("AF37").Value = 80
$ws.Range("AH37").Value = 15
$ws.Range("AJ37").Value = 40
$ws.Range("AN37").Value = 4.05
$ws.Range("AO37").Value = 11.75
$ws.Range("AP37").Value = 22
$ws.Range("AU37").Value = 7.3
$ws.Range("BB37").Value = 350
$ws.Range("G38").Value = 10.25
$ws.Range("H38").Value = 5.1
$ws.Range("J38").Value = 8.75
$ws.Range("K38").Value = 2.47
$ws.Range("N38").Value = 8.5
$ws.Range("O38").Value = 1.22
$ws.Range("U38").Value = 2.3
$ws.Range("W38").Value = 25
$ws.Range("X38").Value = 80
$ws.Range("Y38").Value = 35
$ws.Range("Z38").Value = 350
$ws.Range("AC38").Value = 8.5
$ws.Range("AE38").Value = 29
$ws.Range("AN38").Value = 10.75
$ws.Range("AO38").Value = 70
$ws.Range("AP38").Value = 65
$ws.Range("AU38").Value = 10
$ws.Range("H39").Value = 4.4
$ws.Range("I39").Value = 6.3
$ws.Range("J39").Value = 1.87
$ws.Range("L39").Value = 5.8
$ws.Range("P39").Value = 4.45
$ws.Range("Q39").Value = 1.53
$ws.Range("R39").Value = 2.35
$ws.Range("S39").Value = 1.29
$ws.Range("T39").Value = 3.3
$ws.Range("U39").Value = 1.72
$ws.Range("X39").Value = 7.7
$ws.Range("AG39").Value = 21
$ws.Range("AI39").Value = 20
$ws.Range("AP39").Value = 14
$ws.Range("AR39").Value = 37
$ws.Range("AT39").Value = 3.3
$ws.Range("AW39").Value = 8
$ws.Range("AX39").Value = 35
$ws.Range("H40").Value = 3.6
$ws.Range("P40").Value = 4.3
$ws.Range("R40").Value = 2.3
$ws.Range("T40").Value = 3.25
$ws.Range("U40").Value = 1.52
$ws.Range("V40").Value = 2.37
$ws.Range("W40").Value = 10.25
$ws.Range("X40").Value = 11.25
$ws.Range("AG40").Value = 14
$ws.Range("AL40").Value = 27
$ws.Range("AQ40").Value = 32
$ws.Range("AT40").Value = 3.25
$ws.Range("AY40").Value = 21
$ws.Range("J41").Value = 4.2
$ws.Range("N41").Value = 7.4
$ws.Range("O41").Value = 1.29
$ws.Range("Q41").Value = 1.87
$ws.Range("R41").Value = 1.87
$ws.Range("T41").Value = 2.72
$ws.Range("U41").Value = 1.75
$ws.Range("V41").Value = 1.98
$ws.Range("W41").Value = 10.75
$ws.Range("X41").Value = 20
$ws.Range("AA41").Value = 35
$ws.Range("AB41").Value = 40
$ws.Range("AC41").Value = 7.4
$ws.Range("AE41").Value = 14.5
$ws.Range("AF41").Value = 65
$ws.Range("AG41").Value = 7.5
$ws.Range("AK41").Value = 14.5
$ws.Range("AM41").Value = 500
$ws.Range("AP41").Value = 28
$ws.Range("AT41").Value = 2.72
$ws.Range("AU41").Value = 7.3
$ws.Range("AX41").Value = 9.5
$ws.Range("BA41").Value = 65
$ws.Range("H42").Value = 3.3
$ws.Range("I42").Value = 3.25
$ws.Range("S42").Value = 1.4
$ws.Range("T42").Value = 2.7
$ws.Range("U42").Value = 1.83
$ws.Range("AA42").Value = 17.5
$ws.Range("AF42").Value = 75
$ws.Range("AG42").Value = 9.5
$ws.Range("AJ42").Value = 45
$ws.Range("AS42").Value = 250
$ws.Range("AT42").Value = 2.7
$ws.Range("AW42").Value = 5.2
$ws.Range("BA42").Value = 120
